# Updates the cryptocurrency price/volume table with latest scraped values.
# For D-column (Price) cells whose new text looks like a plain number (e.g. "226.79"),
# Excel would otherwise auto-convert the cell to a numeric type and mangle the display
# (trailing zeros lost, floating point noise, etc.). To keep them as text - exactly like
# the original inline strings - we temporarily force a text NumberFormat, assign the
# value, then restore the cell style to "Normal" so no visible formatting changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '34.385.98'
$ws.Range('E2').Value = '  +0.84%  '
# Row 3
$ws.Range('D3').Value = '1.795.72'
$ws.Range('E3').Value = '  +0.51%  '
# Row 4
$ws.Range('E4').Value = '  +0.12%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.05%  '
# Row 6
$ws.Range('E6').Value = '  +1.41%  '
# Row 7
$ws.Range('E7').Value = '  +0.08%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.55'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.10%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.295'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.17%  '
# Row 10
$ws.Range('E10').Value = '  +0.25%  '
# Row 11
$ws.Range('E11').Value = '  +0.67%  '
# Row 12
$ws.Range('D12').Value = '2.055.05'
$ws.Range('E12').Value = '  +0.55%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.10'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.03%  '
# Row 14
$ws.Range('D14').Value = '1.781.01'
$ws.Range('E14').Value = '  -0.22%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.630'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.72%  '
# Row 16
$ws.Range('D16').Value = '34.377.09'
$ws.Range('E16').Value = '  +0.97%  '
# Row 17
$ws.Range('E17').Value = '  +0.80%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.32'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.49%  '
# Row 19
$ws.Range('D19').Value = '0.0₃0801'
$ws.Range('E19').Value = '  +2.93%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '246.73'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.39%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.99'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.54%  '
# Row 22
$ws.Range('E22').Value = '  +0.07%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.14'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.16%  '
# Row 24
$ws.Range('E24').Value = '  +0.42%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.73'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.65%  '
# Row 26
$ws.Range('E26').Value = '  +0.56%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.40'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.54%  '
# Row 28
$ws.Range('E28').Value = '  +2.02%  '
# Row 29
$ws.Range('E29').Value = '  +0.04%  '
# Row 30
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.02%  '
# Row 31
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0521'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.79%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.90'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.19%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.78'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.40%  '
# Row 34
$ws.Range('E34').Value = '  +1.46%  '
# Row 35
$ws.Range('D35').Value = '1.441.55'
$ws.Range('E35').Value = '  -0.49%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.64'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.39%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.667'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.25%  '
# Row 38
$ws.Range('E38').Value = '  +2.05%  '
# Row 39
$ws.Range('E39').Value = '  -0.82%  '
# Row 40
$ws.Range('E40').Value = '  +4.70%  '
# Row 42
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.76'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.94%  '
# Row 43
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.933'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.67%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.86'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.50%  '
# Row 45
$ws.Range('E45').Value = '  +3.33%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.08'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.69%  '
# Row 47
$ws.Range('E47').Value = '  -0.05%  '
# Row 48
$ws.Range('D48').Value = '1.950.67'
$ws.Range('E48').Value = '  +0.27%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.55'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.64%  '
# Row 50
$ws.Range('E50').Value = '  +0.06%  '
# Row 51
$ws.Range('D51').Value = '0.0₆0128'
$ws.Range('E51').Value = '  -6.26%  '
